# Add a new commemorative coin entry (Saarland, Saarschleife 2025) to the
# "2€" worksheet, replicating rows 35 and 36 after the existing 2024 entry
# (rows 33/34), together with its conditional formatting and number of
# mintage figures, then point the active selection at I37.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Bring over the formatting (styles) used by the previous pair of
#        rows (33/34) so the new rows 35/36 look the same. ---------------
$ws.Range("A33:Q33").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A34:Q34").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# --- 2. Row 35 values -----------------------------------------------------
$ws.Range("A35").Value = 2025
$ws.Range("B35").Value = "Saarland, Saarschleife"
$ws.Range("C35").Value = "Federal states of Germany"
$ws.Range("D35").Value = "Obv: With mint letter"
$ws.Range("E35").Value = "Rev: new map of Europe"
$ws.Range("F35").Value = "6.116.000"
$ws.Range("G35").Value = "6.416.000"
$ws.Range("H35").Value = "7.316.000"
$ws.Range("I35").Value = "4.316.000"
$ws.Range("J35").Value = "6.416.000"
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("N35").Value = 0
$ws.Range("O35").Value = 0

# --- 3. Row 36 values (second mint-letter row, most cells blank) ---------
$ws.Range("A36").Value = 2025
$ws.Range("D36").Value = "Obv: With mint letter"
$ws.Range("E36").Value = "Rev: new map of Europe"
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 0
$ws.Range("O36").Value = 0

# --- 4. "Can exchange" helper formulas in column P ------------------------
$ws.Range("P35:P36").FormulaR1C1 = "=IF(OR(AND(RC[-5]>1,RC[-5]<>""-""),AND(RC[-4]>1,RC[-4]<>""-""),AND(RC[-3]>1,RC[-3]<>""-""),AND(RC[-2]>1,RC[-2]<>""-""),AND(RC[-1]>1,RC[-1]<>""-"")),""Can exchange"","""")"

# --- 5. Conditional formatting for the new rows (same rules as the rest
#        of the table: highlight "-" values, plus a red/yellow/green
#        color scale) ------------------------------------------------------
function Add-MintageConditionalFormatting($range, $cellRef) {
    $containsText = $range.FormatConditions.Add(9, 8, '=NOT(ISERROR(SEARCH(("*-"),(' + $cellRef + '))))')
    $containsText.Formula1 = '=NOT(ISERROR(SEARCH(("*-"),(' + $cellRef + '))))'
    $containsText.Text = "*-"
    $containsText.Interior.Pattern = 1
    $containsText.Interior.Color = 16770459

    $colorScale = $range.FormatConditions.AddColorScale(3)
    $crit = $colorScale.ColorScaleCriteria

    $c1 = $crit.Item(1)
    $c1.Type = 4
    $c1.Value = 0
    $c1.FormatColor.Color = 10461183

    $c2 = $crit.Item(2)
    $c2.Type = 4
    $c2.Value = 1
    $c2.FormatColor.Color = 11722961

    $c3 = $crit.Item(3)
    $c3.Type = 4
    $c3.Value = 10
    $c3.FormatColor.Color = 5287936
}

Add-MintageConditionalFormatting $ws.Range("K35:O35") "K35"
Add-MintageConditionalFormatting $ws.Range("K36:O36") "K36"

# --- 6. Move the selection on the frozen-pane view to I37, as in the
#        edited workbook. ---------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("I37").Select() | Out-Null
